# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Tue Apr  4 21:53:48 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "28.263.92"
$ws.Cells.Item(2, 5).Value = "  +1.99%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.874.63"
$ws.Cells.Item(3, 5).Value = "  +4.42%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.41%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "312.17"
$ws.Cells.Item(5, 5).Value = "  +2.04%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "1.002"
$ws.Cells.Item(6, 5).Value = "  +0.17%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.5062"
$ws.Cells.Item(7, 5).Value = "  +2.01%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.3940"
$ws.Cells.Item(8, 5).Value = "  +2.32%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.09624"
$ws.Cells.Item(9, 5).Value = "  +2.22%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "1.146"
$ws.Cells.Item(10, 5).Value = "  +5.13%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "40.87"
$ws.Cells.Item(11, 5).Value = "  +1.13%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "6.486"

# Row 13
$ws.Cells.Item(13, 4).Value = "21.02"
$ws.Cells.Item(13, 5).Value = "  +3.07%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "1.882.84"
$ws.Cells.Item(14, 5).Value = "  +5.03%  "

# Row 15
$ws.Cells.Item(15, 2).Value = "Chainlink"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(15, 4).Value = "7.435"
$ws.Cells.Item(15, 5).Value = "  +4.33%  "

# Row 16
$ws.Cells.Item(16, 2).Value = "BinanceUSD"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(16, 4).Value = "1.002"
$ws.Cells.Item(16, 5).Value = "  +0.42%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +2.55%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'93.00"
$ws.Cells.Item(18, 5).Value = "  +1.20%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.93%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "17.61"
$ws.Cells.Item(20, 5).Value = "  +3.63%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +0.13%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "6.183"
$ws.Cells.Item(22, 5).Value = "  +5.03%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "28.330.85"
$ws.Cells.Item(23, 5).Value = "  +2.22%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "11.31"
$ws.Cells.Item(24, 5).Value = "  +3.49%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "2.306"
$ws.Cells.Item(25, 5).Value = "  +3.79%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +7.60%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "2.097.19"
$ws.Cells.Item(27, 5).Value = "  +4.94%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "21.24"
$ws.Cells.Item(28, 5).Value = "  +4.31%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'158.90"
$ws.Cells.Item(29, 5).Value = "  +1.35%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "127.62"
$ws.Cells.Item(30, 5).Value = "  +1.30%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'0.1070"
$ws.Cells.Item(31, 5).Value = "  +0.34%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "1.069"
$ws.Cells.Item(32, 5).Value = "  +1.85%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +2.62%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "3.624"
$ws.Cells.Item(34, 5).Value = "  +0.60%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "9.554"
$ws.Cells.Item(35, 5).Value = "  +8.14%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'0.06710"
$ws.Cells.Item(36, 5).Value = "  -1.13%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "0.02392"
$ws.Cells.Item(37, 5).Value = "  +4.27%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "0.2194"
$ws.Cells.Item(38, 5).Value = "  +3.43%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "TheSandbox"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(39, 4).Value = "0.6382"
$ws.Cells.Item(39, 5).Value = "  +4.51%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Aptos"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(40, 4).Value = "11.52"
$ws.Cells.Item(40, 5).Value = "  +1.56%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +2.14%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "1.188"
$ws.Cells.Item(42, 5).Value = "  +4.52%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +0.12%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "13.47"
$ws.Cells.Item(44, 5).Value = "  +4.50%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "0.5994"
$ws.Cells.Item(45, 5).Value = "  +2.58%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "3.662"
$ws.Cells.Item(46, 5).Value = "  -0.02%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "1.268"
$ws.Cells.Item(47, 5).Value = "  +0.40%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "2.007"
$ws.Cells.Item(48, 5).Value = "  +4.31%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "124.22"
$ws.Cells.Item(49, 5).Value = "  +0.97%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "1.196"
$ws.Cells.Item(50, 5).Value = "  +2.57%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "0.06858"
$ws.Cells.Item(51, 5).Value = "  +2.68%  "
